$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Insert two new columns: D (RF "2 att") and G (KNN "2 att")
#    This shifts the old KNN columns D,E -> E,F and carries their
#    formatting along automatically.
# ---------------------------------------------------------------
$ws.Columns.Item(4).Insert()   # new column D (between old C and D)
$ws.Columns.Item(7).Insert()   # new column G (after the shifted KNN block)

# ---------------------------------------------------------------
# 2) Re-create the header merges across the now-widened blocks
# ---------------------------------------------------------------
$ws.Range("B2:D2").Merge()
$ws.Range("E2:G2").Merge()

# ---------------------------------------------------------------
# 3) Fill in the new "2 att" columns with the new model results
#    (RF -> column D, KNN -> column G)
# ---------------------------------------------------------------
$ws.Range("D3").Value = "2 att"
$ws.Range("D4").Value = "(mtry:2) rmse:153.41, Rsq:0.94, mae:84.96"
$ws.Range("D5").Value = "rmse:249.41, Rsq:0.87, mae:148.6"

$ws.Range("G3").Value = "2 att"
$ws.Range("G4").Value = "(k:5) rmse:241.47, Rsq:0.85, mae:115.35"
$ws.Range("G5").Value = "rmse:312.05, Rsq:0.82, mae:179.98"

# ---------------------------------------------------------------
# 4) Column widths
# ---------------------------------------------------------------
$ws.Range("B1").ColumnWidth = 30.833333333333336
$ws.Range("C1").ColumnWidth = 30.833333333333336
$ws.Range("D1").ColumnWidth = 33.166666666666664
$ws.Range("E1").ColumnWidth = 27.833333333333336
$ws.Range("F1").ColumnWidth = 27.833333333333336
$ws.Range("G1").ColumnWidth = 27.666666666666668

# ---------------------------------------------------------------
# 5) Border: thin right-hand border on column D, separating the
#    RF block from the KNN block
# ---------------------------------------------------------------
$ws.Range("D2:D5").Borders.Item(10).Weight = 2

# ---------------------------------------------------------------
# 6) Alignment: header row 2 stays centered across the full width
# ---------------------------------------------------------------
$ws.Range("B2:G2").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# 7) Sheet view: zoom 90%, scrolled so column B is the leftmost
#    visible column, selection moved to E15
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Application.ActiveWindow.ScrollColumn = 2
$null = $ws.Range("E15").Select()

Write-Host "done"
